$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-06-01 Saturday" "2024-06-02 Sunday"

Replace-Text "45×71=" "43×91="
Replace-Text "82×95=" "88×66="
Replace-Text "92×48=" "67×27="
Replace-Text "25×48=" "47×48="
Replace-Text "51×51=" "91×91="
Replace-Text "68×82=" "28×49="
Replace-Text "73×35=" "52×70="
Replace-Text "51×13=" "90×95="
Replace-Text "20×20=" "88×77="
Replace-Text "93×51=" "99×48="
Replace-Text "71×20=" "21×62="
Replace-Text "15×77=" "88×54="
Replace-Text "31×19=" "63×38="
Replace-Text "84×36=" "30×41="
Replace-Text "35×45=" "23×61="
Replace-Text "23×63=" "28×37="
Replace-Text "47×19=" "34×41="
Replace-Text "49×77=" "51×25="
Replace-Text "91×21=" "31×57="
Replace-Text "11×33=" "75×91="
Replace-Text "62×29=" "73×40="
Replace-Text "30×90=" "63×74="
Replace-Text "31×66=" "29×76="
Replace-Text "77×65=" "62×71="
Replace-Text "69×56=" "57×62="
